# repull data, push all data, mean calculation
# Update the dSF column (column F) values on the active sheet to reflect
# the repulled/updated source data for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    5  = -4
    17 = -3
    18 = -2
    20 = 4
    22 = -8
    23 = 2
    26 = -1
    31 = -2
    32 = -2
    34 = 3
    35 = -3
    36 = 0
    40 = -4
    45 = 0
    49 = -1
    53 = -1
    62 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
